$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.888.91"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.144.27"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.135.36"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "3.666.67"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.742.54"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.140.46"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.84%  "
$ws.Range("E29").Value = "  +8.12%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").Value = "0.0₃0842"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "454.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.293"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.22%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "2.913.83"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.25%  "
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.49%  "
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("E51").Value = "  -1.07%  "
